$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analytics")

# Update the event name in C5 ("event" -> "event1")
$ws.Range("C5").Value = "event1"

# Update the active selection to C5 (as reflected in the saved view state)
$ws.Range("C5").Select()
